$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 199.75
$ws.Range("I5").Value = 265.66666
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 265.66666
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = -150.66666
$ws.Range("N5").Value = -232
$ws.Range("H15").Value = 2360.4285
$ws.Range("I15").Value = 2360.4285
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 7081.2855
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -6912.2855
$ws.Range("H19").Value = 2717.4
$ws.Range("I19").Value = 2749.75
$ws.Range("J19").Value = 2709.3125
$ws.Range("K19").Value = 2749.75
$ws.Range("L19").Value = 2709.3125
$ws.Range("M19").Value = -2574.75
$ws.Range("N19").Value = -3059.3125
$ws.Range("H34").Value = 11788.6
$ws.Range("I34").Value = 4736
$ws.Range("J34").Value = 39999
$ws.Range("K34").Value = 4736
$ws.Range("L34").Value = 39999
$ws.Range("M34").Value = -4533
$ws.Range("H36").Value = 11788.6
$ws.Range("I36").Value = 4736
$ws.Range("J36").Value = 39999
$ws.Range("K36").Value = 4736
$ws.Range("L36").Value = 39999
$ws.Range("M36").Value = -4021
$ws.Range("H55").Value = 80
$ws.Range("I55").Value = 44.666668
$ws.Range("J55").Value = 93.25
$ws.Range("K55").Value = 44.666668
$ws.Range("L55").Value = 93.25
$ws.Range("M55").Value = 169.333332
$ws.Range("N55").Value = -521.25
$ws.Range("H80").Value = 799.7778
$ws.Range("I80").Value = 466.66666
$ws.Range("J80").Value = 966.3333
$ws.Range("K80").Value = 1399.99998
$ws.Range("L80").Value = 2898.9999
$ws.Range("M80").Value = -401.9999800000001
$ws.Range("N80").Value = -4894.9999
$ws.Range("H83").Value = 799.7778
$ws.Range("I83").Value = 466.66666
$ws.Range("J83").Value = 966.3333
$ws.Range("K83").Value = 4199.99994
$ws.Range("L83").Value = 8696.9997
$ws.Range("M83").Value = 792.0000600000003
$ws.Range("N83").Value = -18680.9997
$ws.Range("H92").Value = 730.3125
$ws.Range("I92").Value = 610.38464
$ws.Range("J92").Value = 1250
$ws.Range("K92").Value = 610.38464
$ws.Range("L92").Value = 1250
$ws.Range("M92").Value = 637.61536
$ws.Range("H98").Value = 5999.2
$ws.Range("I98").Value = 999.3333
$ws.Range("J98").Value = 13499
$ws.Range("K98").Value = 999.3333
$ws.Range("L98").Value = 13499
$ws.Range("M98").Value = 498.6667
$ws.Range("H106").Value = 7294.375
$ws.Range("I106").Value = 7294.375
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 7294.375
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -6663.375
$ws.Range("H111").Value = 675.8
$ws.Range("I111").Value = 652.2
$ws.Range("J111").Value = 699.4
$ws.Range("K111").Value = 1956.6
$ws.Range("L111").Value = 2098.2
$ws.Range("M111").Value = 1110.4
$ws.Range("H116").Value = 6406.375
$ws.Range("I116").Value = 6249.5
$ws.Range("J116").Value = 6458.6665
$ws.Range("K116").Value = 6249.5
$ws.Range("L116").Value = 6458.6665
$ws.Range("M116").Value = -2807.5
$ws.Range("N116").Value = -13342.6665
$ws.Range("H122").Value = 5999.2
$ws.Range("I122").Value = 999.3333
$ws.Range("J122").Value = 13499
$ws.Range("K122").Value = 2997.9999
$ws.Range("L122").Value = 40497
$ws.Range("M122").Value = -547.9998999999998
$ws.Range("H125").Value = 17331
$ws.Range("I125").Value = 20998
$ws.Range("J125").Value = 9997
$ws.Range("K125").Value = 188982
$ws.Range("L125").Value = 89973
$ws.Range("M125").Value = -186522
$ws.Range("H132").Value = 3383.9412
$ws.Range("I132").Value = 3014.9768
$ws.Range("J132").Value = 5367.125
$ws.Range("K132").Value = 9044.930399999999
$ws.Range("L132").Value = 16101.375
$ws.Range("M132").Value = -6514.930399999999
$ws.Range("H135").Value = 2515.2727
$ws.Range("I135").Value = 2728.9
$ws.Range("J135").Value = 379
$ws.Range("K135").Value = 24560.1
$ws.Range("L135").Value = 3411
$ws.Range("M135").Value = -22025.1
$ws.Range("N135").Value = -8481
$ws.Range("H137").Value = 1622.1464
$ws.Range("I137").Value = 1298.7646
$ws.Range("J137").Value = 3192.8572
$ws.Range("K137").Value = 3896.2938
$ws.Range("L137").Value = 9578.571599999999
$ws.Range("M137").Value = -1346.2938
$ws.Range("H138").Value = 3146.4167
$ws.Range("I138").Value = 2203.75
$ws.Range("J138").Value = 3489.2046
$ws.Range("K138").Value = 6611.25
$ws.Range("L138").Value = 10467.6138
$ws.Range("M138").Value = -1471.25
$ws.Range("N138").Value = -20747.6138
$ws.Range("H141").Value = 3599
$ws.Range("I141").Value = 3776.5557
$ws.Range("J141").Value = 2800
$ws.Range("K141").Value = 11329.6671
$ws.Range("L141").Value = 8400
$ws.Range("M141").Value = -6149.667099999999
$ws.Range("N141").Value = -18760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 5495
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 5495
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5495
$ws.Range("N17").Value = -5841
$ws.Range("H32").Value = 6418.2925
$ws.Range("I32").Value = 1419.7142
$ws.Range("J32").Value = 17184.46
$ws.Range("K32").Value = 1419.7142
$ws.Range("L32").Value = 17184.46
$ws.Range("M32").Value = -1132.7142
$ws.Range("N32").Value = -17758.46
$ws.Range("H61").Value = 8037.077
$ws.Range("I61").Value = 5862.273
$ws.Range("J61").Value = 19998.5
$ws.Range("K61").Value = 5862.273
$ws.Range("L61").Value = 19998.5
$ws.Range("M61").Value = -5650.273
$ws.Range("H74").Value = 5642.778
$ws.Range("I74").Value = 3969.4285
$ws.Range("J74").Value = 11499.5
$ws.Range("K74").Value = 3969.4285
$ws.Range("L74").Value = 11499.5
$ws.Range("M74").Value = -3095.4285
$ws.Range("H77").Value = 5642.778
$ws.Range("I77").Value = 3969.4285
$ws.Range("J77").Value = 11499.5
$ws.Range("K77").Value = 19847.1425
$ws.Range("L77").Value = 57497.5
$ws.Range("M77").Value = -15479.1425
$ws.Range("H114").Value = 64250
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 64250
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 64250
$ws.Range("N114").Value = -72928
$ws.Range("H122").Value = 5156.409
$ws.Range("I122").Value = 4658.0835
$ws.Range("J122").Value = 7398.875
$ws.Range("K122").Value = 13974.2505
$ws.Range("L122").Value = 22196.625
$ws.Range("M122").Value = -11524.2505
$ws.Range("N122").Value = -27096.625
$ws.Range("H132").Value = 2394.64
$ws.Range("I132").Value = 2202.75
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 6608.25
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -4078.25
$ws.Range("H136").Value = 8037.077
$ws.Range("I136").Value = 5862.273
$ws.Range("J136").Value = 19998.5
$ws.Range("K136").Value = 17586.819
$ws.Range("L136").Value = 59995.5
$ws.Range("M136").Value = -15036.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = ""
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = ""
$ws.Range("H64").Value = 365.66666
$ws.Range("I64").Value = 298.5
$ws.Range("J64").Value = 500
$ws.Range("K64").Value = 298.5
$ws.Range("L64").Value = 500
$ws.Range("M64").Value = -73.5
$ws.Range("N64").Value = -950
$ws.Range("H67").Value = 365.66666
$ws.Range("I67").Value = 298.5
$ws.Range("J67").Value = 500
$ws.Range("K67").Value = 298.5
$ws.Range("L67").Value = 500
$ws.Range("M67").Value = 481.5
$ws.Range("N67").Value = -2060
$ws.Range("H86").Value = 1289.8422
$ws.Range("I86").Value = 861.2
$ws.Range("J86").Value = 2897.25
$ws.Range("K86").Value = 861.2
$ws.Range("L86").Value = 2897.25
$ws.Range("M86").Value = 261.8
$ws.Range("H89").Value = 1289.8422
$ws.Range("I89").Value = 861.2
$ws.Range("J89").Value = 2897.25
$ws.Range("K89").Value = 4306
$ws.Range("L89").Value = 14486.25
$ws.Range("M89").Value = 1310
$ws.Range("H99").Value = 2666.2173
$ws.Range("I99").Value = 2666.2173
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2666.2173
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1168.2173
$ws.Range("H105").Value = 2637.6155
$ws.Range("I105").Value = 2637.6155
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2637.6155
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -890.6154999999999
$ws.Range("H134").Value = 3892.2942
$ws.Range("I134").Value = 3891.8125
$ws.Range("J134").Value = 3900
$ws.Range("K134").Value = 11675.4375
$ws.Range("L134").Value = 11700
$ws.Range("M134").Value = -9140.4375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 15449
$ws.Range("I16").Value = 1123.875
$ws.Range("J16").Value = 72749.5
$ws.Range("K16").Value = 1123.875
$ws.Range("L16").Value = 72749.5
$ws.Range("M16").Value = -836.875
$ws.Range("N16").Value = -73323.5
$ws.Range("H22").Value = 1062
$ws.Range("I22").Value = 1023.5
$ws.Range("J22").Value = 1241.6666
$ws.Range("K22").Value = 1023.5
$ws.Range("L22").Value = 1241.6666
$ws.Range("M22").Value = -673.5
$ws.Range("H31").Value = 6560.431
$ws.Range("I31").Value = 6341.3687
$ws.Range("J31").Value = 6690.5
$ws.Range("K31").Value = 6341.3687
$ws.Range("L31").Value = 6690.5
$ws.Range("M31").Value = -6046.3687
$ws.Range("H34").Value = 6560.431
$ws.Range("I34").Value = 6341.3687
$ws.Range("J34").Value = 6690.5
$ws.Range("K34").Value = 6341.3687
$ws.Range("L34").Value = 6690.5
$ws.Range("M34").Value = -6139.3687
$ws.Range("H96").Value = 10027.6
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 10027.6
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 10027.6
$ws.Range("N96").Value = -15519.6
$ws.Range("H97").Value = 97197
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 97197
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 97197
$ws.Range("N97").Value = -99179
$ws.Range("H99").Value = 5231.8667
$ws.Range("I99").Value = 4832.1113
$ws.Range("J99").Value = 5831.5
$ws.Range("K99").Value = 4832.1113
$ws.Range("L99").Value = 5831.5
$ws.Range("M99").Value = -3334.1113
$ws.Range("N99").Value = -8827.5
$ws.Range("H113").Value = 15449
$ws.Range("I113").Value = 1123.875
$ws.Range("J113").Value = 72749.5
$ws.Range("K113").Value = 1123.875
$ws.Range("L113").Value = 72749.5
$ws.Range("M113").Value = 1046.125
$ws.Range("N113").Value = -77089.5
$ws.Range("H115").Value = 36900
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 36900
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 36900
$ws.Range("N115").Value = -39250
$ws.Range("H122").Value = 103837.63
$ws.Range("I122").Value = 134107.78
$ws.Range("J122").Value = 4378.5713
$ws.Range("K122").Value = 402323.34
$ws.Range("L122").Value = 13135.7139
$ws.Range("M122").Value = -399873.34
$ws.Range("H126").Value = 5231.8667
$ws.Range("I126").Value = 4832.1113
$ws.Range("J126").Value = 5831.5
$ws.Range("K126").Value = 14496.3339
$ws.Range("L126").Value = 17494.5
$ws.Range("M126").Value = -12026.3339
$ws.Range("N126").Value = -22434.5
$ws.Range("H132").Value = 1253.2
$ws.Range("I132").Value = 944.3333
$ws.Range("J132").Value = 2488.6667
$ws.Range("K132").Value = 2832.9999
$ws.Range("L132").Value = 7466.000100000001
$ws.Range("M132").Value = -302.9998999999998
$ws.Range("H134").Value = 1055.8182
$ws.Range("I134").Value = 891.0625
$ws.Range("J134").Value = 1495.1666
$ws.Range("K134").Value = 2673.1875
$ws.Range("L134").Value = 4485.4998
$ws.Range("M134").Value = -138.1875
$ws.Range("H141").Value = 215832
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 215832
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 215832
$ws.Range("N141").Value = -226192

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 106895096
$ws.Range("I4").Value = 143985860
$ws.Range("J4").Value = 20350000
$ws.Range("K4").Value = 431957580
$ws.Range("L4").Value = 61050000
$ws.Range("M4").Value = -431957468
$ws.Range("N4").Value = -61050224
$ws.Range("H5").Value = 499.82352
$ws.Range("I5").Value = 385.7857
$ws.Range("J5").Value = 1032
$ws.Range("K5").Value = 1157.3571
$ws.Range("L5").Value = 3096
$ws.Range("M5").Value = -1045.3571
$ws.Range("N5").Value = -3320
$ws.Range("H92").Value = 716.1818
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 716.1818
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2148.5454
$ws.Range("N92").Value = -4644.5454
$ws.Range("H107").Value = 458.47058
$ws.Range("I107").Value = 524.8333
$ws.Range("J107").Value = 422.27274
$ws.Range("K107").Value = 1574.4999
$ws.Range("L107").Value = 1266.81822
$ws.Range("M107").Value = 345.5001
$ws.Range("N107").Value = -5106.81822
$ws.Range("H122").Value = 2451.2083
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 2527.348
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 22746.132
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -27646.132
$ws.Range("H129").Value = 1157.5
$ws.Range("I129").Value = 889
$ws.Range("J129").Value = 2500
$ws.Range("K129").Value = 2667
$ws.Range("L129").Value = 7500
$ws.Range("M129").Value = 2333
$ws.Range("N129").Value = -17500
$ws.Range("H132").Value = 558.1667
$ws.Range("I132").Value = 558.1667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5023.5003
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2493.5003
$ws.Range("H135").Value = 499.82352
$ws.Range("I135").Value = 385.7857
$ws.Range("J135").Value = 1032
$ws.Range("K135").Value = 3472.0713
$ws.Range("L135").Value = 9288
$ws.Range("M135").Value = -937.0713000000001
$ws.Range("N135").Value = -14358
$ws.Range("H136").Value = 5699.5454
$ws.Range("I136").Value = 4654.1816
$ws.Range("J136").Value = 7790.273
$ws.Range("K136").Value = 13962.5448
$ws.Range("L136").Value = 23370.819
$ws.Range("M136").Value = -8862.5448
$ws.Range("H139").Value = 3364.05
$ws.Range("I139").Value = 3277.9473
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 9833.841899999999
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = -4693.841899999999
$ws.Range("N139").Value = -25280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6188.778
$ws.Range("I80").Value = 2003
$ws.Range("J80").Value = 8281.666999999999
$ws.Range("K80").Value = 2003
$ws.Range("L80").Value = 8281.666999999999
$ws.Range("M80").Value = -1005
$ws.Range("N80").Value = -10277.667
$ws.Range("H83").Value = 6188.778
$ws.Range("I83").Value = 2003
$ws.Range("J83").Value = 8281.666999999999
$ws.Range("K83").Value = 10015
$ws.Range("L83").Value = 41408.335
$ws.Range("M83").Value = -5023
$ws.Range("N83").Value = -51392.335
$ws.Range("H102").Value = 1863.2273
$ws.Range("I102").Value = 1388.3889
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 1388.3889
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = 233.6111000000001
$ws.Range("N102").Value = -7244
$ws.Range("H107").Value = 834.53845
$ws.Range("I107").Value = 785.2
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 785.2
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 1134.8
$ws.Range("H122").Value = 8671.333000000001
$ws.Range("I122").Value = 8007
$ws.Range("J122").Value = 9003.5
$ws.Range("K122").Value = 24021
$ws.Range("L122").Value = 27010.5
$ws.Range("M122").Value = -21571
$ws.Range("N122").Value = -31910.5
$ws.Range("H126").Value = 9304.333000000001
$ws.Range("I126").Value = 7441.8887
$ws.Range("J126").Value = 11166.777
$ws.Range("K126").Value = 22325.6661
$ws.Range("L126").Value = 33500.331
$ws.Range("M126").Value = -19855.6661
$ws.Range("N126").Value = -38440.331
$ws.Range("H132").Value = 4515.6665
$ws.Range("I132").Value = 4469.5483
$ws.Range("J132").Value = 4801.6
$ws.Range("K132").Value = 13408.6449
$ws.Range("L132").Value = 14404.8
$ws.Range("M132").Value = -10878.6449
$ws.Range("N132").Value = -19464.8
$ws.Range("H133").Value = 89302.75
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 89302.75
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 89302.75
$ws.Range("N133").Value = -99422.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3876.6924
$ws.Range("I7").Value = 3566.3333
$ws.Range("J7").Value = 4142.7144
$ws.Range("K7").Value = 3566.3333
$ws.Range("L7").Value = 4142.7144
$ws.Range("M7").Value = -3454.3333
$ws.Range("H22").Value = 1949.5714
$ws.Range("I22").Value = 1559.8
$ws.Range("J22").Value = 2924
$ws.Range("K22").Value = 1559.8
$ws.Range("L22").Value = 2924
$ws.Range("M22").Value = -1264.8
$ws.Range("N22").Value = -3514
$ws.Range("H27").Value = 1949.5714
$ws.Range("I27").Value = 1559.8
$ws.Range("J27").Value = 2924
$ws.Range("K27").Value = 1559.8
$ws.Range("L27").Value = 2924
$ws.Range("M27").Value = -1452.8
$ws.Range("N27").Value = -3138
$ws.Range("H61").Value = 2731.818
$ws.Range("I61").Value = 2206.5715
$ws.Range("J61").Value = 3651
$ws.Range("K61").Value = 2206.5715
$ws.Range("L61").Value = 3651
$ws.Range("M61").Value = -2004.5715
$ws.Range("H68").Value = 2373.5833
$ws.Range("I68").Value = 1323.625
$ws.Range("J68").Value = 4473.5
$ws.Range("K68").Value = 1323.625
$ws.Range("L68").Value = 4473.5
$ws.Range("M68").Value = -574.625
$ws.Range("N68").Value = -5971.5
$ws.Range("H71").Value = 2373.5833
$ws.Range("I71").Value = 1323.625
$ws.Range("J71").Value = 4473.5
$ws.Range("K71").Value = 6618.125
$ws.Range("L71").Value = 22367.5
$ws.Range("M71").Value = -2874.125
$ws.Range("N71").Value = -29855.5
$ws.Range("H82").Value = 1075
$ws.Range("I82").Value = 1020.44446
$ws.Range("J82").Value = 1129.5555
$ws.Range("K82").Value = 1020.44446
$ws.Range("L82").Value = 1129.5555
$ws.Range("M82").Value = -659.44446
$ws.Range("N82").Value = -1851.5555
$ws.Range("H85").Value = 1075
$ws.Range("I85").Value = 1020.44446
$ws.Range("J85").Value = 1129.5555
$ws.Range("K85").Value = 1020.44446
$ws.Range("L85").Value = 1129.5555
$ws.Range("M85").Value = 227.55554
$ws.Range("N85").Value = -3625.5555
$ws.Range("H93").Value = 1389.5333
$ws.Range("I93").Value = 1330.3636
$ws.Range("J93").Value = 1552.25
$ws.Range("K93").Value = 1330.3636
$ws.Range("L93").Value = 1552.25
$ws.Range("M93").Value = -82.36359999999991
$ws.Range("N93").Value = -4048.25
$ws.Range("H113").Value = 2731.818
$ws.Range("I113").Value = 2206.5715
$ws.Range("J113").Value = 3651
$ws.Range("K113").Value = 2206.5715
$ws.Range("L113").Value = 3651
$ws.Range("M113").Value = -36.57150000000001
$ws.Range("H126").Value = 3876.6924
$ws.Range("I126").Value = 3566.3333
$ws.Range("J126").Value = 4142.7144
$ws.Range("K126").Value = 10698.9999
$ws.Range("L126").Value = 12428.1432
$ws.Range("M126").Value = -8228.999899999999
$ws.Range("H132").Value = 7955.353
$ws.Range("I132").Value = 7955.353
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 23866.059
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -21336.059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 7999.8
$ws.Range("I33").Value = 3000
$ws.Range("J33").Value = 27999
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 27999
$ws.Range("M33").Value = -2750
$ws.Range("N33").Value = -28499
$ws.Range("H36").Value = 7999.8
$ws.Range("I36").Value = 3000
$ws.Range("J36").Value = 27999
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 27999
$ws.Range("M36").Value = -2750
$ws.Range("N36").Value = -28499
$ws.Range("H88").Value = 46600
$ws.Range("I88").Value = 20000
$ws.Range("J88").Value = 59900
$ws.Range("K88").Value = 20000
$ws.Range("L88").Value = 59900
$ws.Range("M88").Value = -19594
$ws.Range("N88").Value = -60712
$ws.Range("H91").Value = 46600
$ws.Range("I91").Value = 20000
$ws.Range("J91").Value = 59900
$ws.Range("K91").Value = 20000
$ws.Range("L91").Value = 59900
$ws.Range("M91").Value = -18596
$ws.Range("N91").Value = -62708
$ws.Range("H107").Value = 140.6842
$ws.Range("I107").Value = 158.13333
$ws.Range("J107").Value = 75.25
$ws.Range("K107").Value = 474.39999
$ws.Range("L107").Value = 225.75
$ws.Range("M107").Value = 1445.60001
$ws.Range("N107").Value = -4065.75
$ws.Range("H122").Value = 4118.5293
$ws.Range("I122").Value = 4225.5386
$ws.Range("J122").Value = 3770.75
$ws.Range("K122").Value = 12676.6158
$ws.Range("L122").Value = 11312.25
$ws.Range("M122").Value = -10226.6158
$ws.Range("H132").Value = 6029.4375
$ws.Range("I132").Value = 4326.737
$ws.Range("J132").Value = 12499.7
$ws.Range("K132").Value = 12980.211
$ws.Range("L132").Value = 37499.10000000001
$ws.Range("M132").Value = -10450.211
$ws.Range("N132").Value = -42559.10000000001
$ws.Range("H136").Value = 7504.4375
$ws.Range("I136").Value = 8472.546
$ws.Range("J136").Value = 5374.6
$ws.Range("K136").Value = 25417.638
$ws.Range("L136").Value = 16123.8
$ws.Range("M136").Value = -22867.638
$ws.Range("N136").Value = -21223.8
$ws.Range("H138").Value = 89332.664
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 89332.664
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 89332.664
$ws.Range("N138").Value = -99612.664
$ws.Range("H139").Value = 110081.78
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 110081.78
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 110081.78
$ws.Range("N139").Value = -120361.78
